$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.089.34"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.015.89"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'226.21"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'54.87"
$ws.Range("E8").Value = "  -5.78%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("D10").Value = "'0.0783"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  -5.20%  "
$ws.Range("D12").Value = "2.313.71"
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "'5.12"
$ws.Range("E16").Value = "  -4.18%  "
$ws.Range("D17").Value = "2.028.98"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "37.009.39"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'68.86"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "0.0₃0816"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "'223.07"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  -6.55%  "
$ws.Range("D26").Value = "'166.22"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("E27").Value = "  -8.08%  "
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.124"
$ws.Range("E29").Value = "  -5.16%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'18.68"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").Value = "'4.42"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").Value = "  -7.69%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("D39").Value = "'5.31"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "1.478.01"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  -5.52%  "
$ws.Range("D42").Value = "'94.79"
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").Value = "'0.0913"
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("D44").Value = "'16.27"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("E45").Value = "  -5.44%  "
$ws.Range("E46").Value = "  -5.95%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.18"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "2.201.44"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").Value = "'44.29"
$ws.Range("E51").Value = "  -3.74%  "
